$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 7779
$ws.Range("K3").Value = 8050
$ws.Range("K4").Value = 1693
$ws.Range("K5").Value = 576
$ws.Range("K6").Value = 8972
$ws.Range("K7").Value = 27070

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K5").Value = 73
$ws.Range("K6").Value = 196
$ws.Range("K8").Value = 1769
$ws.Range("K10").Value = 163
$ws.Range("K14").Value = 128
$ws.Range("K15").Value = 276
$ws.Range("K18").Value = 184
$ws.Range("K19").Value = 778
$ws.Range("K20").Value = 666
$ws.Range("K24").Value = 87
$ws.Range("K25").Value = 131
$ws.Range("K29").Value = 1495
$ws.Range("K33").Value = 1136
$ws.Range("K34").Value = 154
$ws.Range("K36").Value = 351
$ws.Range("K37").Value = 890
$ws.Range("K40").Value = 62
$ws.Range("K42").Value = 1003
$ws.Range("K47").Value = 188
$ws.Range("K48").Value = 340
$ws.Range("K49").Value = 151
$ws.Range("K51").Value = 350
$ws.Range("K52").Value = 694
$ws.Range("K53").Value = 343
$ws.Range("K54").Value = 530
$ws.Range("K55").Value = 297
$ws.Range("K63").Value = 76
$ws.Range("K64").Value = 164
$ws.Range("K67").Value = 1053
$ws.Range("K70").Value = 49
$ws.Range("K74").Value = 28
$ws.Range("K75").Value = 88
$ws.Range("K83").Value = 577
$ws.Range("K84").Value = 219
$ws.Range("K85").Value = 1246
$ws.Range("K91").Value = 325
$ws.Range("K99").Value = 455
$ws.Range("K101").Value = 27070

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 413
$ws.Range("K3").Value = 431
$ws.Range("K7").Value = 1246

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 192
$ws.Range("K7").Value = 694

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 92
$ws.Range("K7").Value = 343

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 494
$ws.Range("K3").Value = 534
$ws.Range("K7").Value = 1769

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 200
$ws.Range("K3").Value = 201
$ws.Range("K7").Value = 577

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 402
$ws.Range("K7").Value = 1136

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 296
$ws.Range("K7").Value = 890

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K2").Value = 121
$ws.Range("K3").Value = 193
$ws.Range("K6").Value = 109
$ws.Range("K7").Value = 455

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 287
$ws.Range("K3").Value = 382
$ws.Range("K6").Value = 298
$ws.Range("K7").Value = 1053

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K3").Value = 89
$ws.Range("K7").Value = 219

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 151

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 121
$ws.Range("K4").Value = 38
$ws.Range("K7").Value = 530

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 420
$ws.Range("K6").Value = 444
$ws.Range("K7").Value = 1495

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K4").Value = 51
$ws.Range("K6").Value = 153
$ws.Range("K7").Value = 340

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 233
$ws.Range("K7").Value = 778

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 77
$ws.Range("K7").Value = 196

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 292
$ws.Range("K6").Value = 384
$ws.Range("K7").Value = 1003

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 163

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K6").Value = 110
$ws.Range("K7").Value = 297

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K5").Value = 10
$ws.Range("K7").Value = 325

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 164

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K6").Value = 194
$ws.Range("K7").Value = 666

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 184

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K3").Value = 112
$ws.Range("K7").Value = 351

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("K3").Value = 41
$ws.Range("K7").Value = 154

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K3").Value = 47
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K3").Value = 56
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 188

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 276

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K2").Value = 44
$ws.Range("K6").Value = 122

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K3").Value = 97
$ws.Range("K7").Value = 350

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K4").Value = 12
$ws.Range("K6").Value = 46

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 28

Write-Host "Applied all 2024-12-25 updates"